# Refresh cryptos list (GitHub Actions scheduled update).
# Note: some Price (column D) values look like plain decimals (e.g. "576.46");
# those are entered with a leading apostrophe so Excel keeps them as text
# (matching the original inline-string cells) instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.874.67'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.101.65'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''576.46'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = '''177.47'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '3.100.42'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').Value = '''36.13'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').Value = '3.620.05'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '66.903.57'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = '''7.04'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '3.102.51'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = '''16.66'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '''480.40'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').Value = '''7.84'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('D24').Value = '''83.70'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').Value = '''12.61'
$ws.Range('E25').Value = '  -3.95%  '
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').Value = '''10.09'
$ws.Range('E27').Value = '  -4.55%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '''7.92'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '''2.29'
$ws.Range('E30').Value = '  -3.09%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').Value = '''27.97'
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('D34').Value = '0.0₃0939'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').Value = '''48.18'
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('E37').Value = '  -4.54%  '
$ws.Range('D38').Value = '''0.942'
$ws.Range('E38').Value = '  -3.38%  '
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('D40').Value = '''49.05'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('D45').Value = '2.801.88'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '''374.14'
$ws.Range('E46').Value = '  -4.15%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '''0.0343'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''135.24'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D50').Value = '''25.45'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('E51').Value = '  +1.72%  '
